# Apply the latest cryptos list values (price + 1h volume change) per cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.176.66"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "2.177.35"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'250.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("E7").Value = "  -7.93%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.577"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.25%  "
$ws.Range("D10").Value = "'59.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "'36.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.14%  "
$ws.Range("E12").Value = "  -4.12%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "'6.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.48%  "
$ws.Range("D15").Value = "2.504.20"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "'14.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "2.194.76"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "41.072.08"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "'71.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").Value = "'228.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").Value = "'2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.90%  "
$ws.Range("D25").Value = "'3.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.29%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'11.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.27%  "
$ws.Range("E28").Value = "  -5.16%  "
$ws.Range("D29").Value = "'3.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.99%  "
$ws.Range("D30").Value = "'168.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("D32").Value = "'20.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").Value = "'5.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.70%  "
$ws.Range("D35").Value = "'0.0757"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("E37").Value = "  -4.35%  "
$ws.Range("D38").Value = "'3.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'24.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.21%  "
$ws.Range("D40").Value = "'0.0307"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("E42").Value = "  +6.11%  "
$ws.Range("E43").Value = "  -8.74%  "
$ws.Range("D44").Value = "'61.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.49%  "
$ws.Range("E45").Value = "  -5.71%  "
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "'0.189"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.69%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("E51").Value = "  -4.08%  "
